# Updating filtered feeds from workflow
# Adds a new row (row 10) to the "Filtered Feeds" sheet for the
# BioCentury article about ODAC voting for Darzalex in smoldering
# multiple myeloma.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$link = "https://www.biocentury.com/article/656014/odac-votes-for-darzalex-in-smoldering-multiple-myeloma"
$keywords = "smoldering multiple myeloma"
$title = "ODAC votes for Darzalex in smoldering multiple myeloma"

# Populate the new row's values.
$ws.Range("A10").Value2 = $link
$ws.Range("B10").Value2 = $keywords
$ws.Range("C10").Value2 = $title

# Turn A10 into a real hyperlink pointing at the article (mirrors the
# other rows, which all have a hyperlink on column A).
$ws.Hyperlinks.Add($ws.Range("A10"), $link) | Out-Null

# Match the "link" column style used by the other data rows (A2:A9).
$ws.Range("A10").Style = $ws.Range("A9").Style
